# Applies the "output generated at 456a3b4" update:
#  - Sheet "展览" (1): refresh "want to go" counts (col F) for several rows,
#    and append a new row 46 for the newly scraped "北京·美漫超级英雄ONLY" event.
#  - Sheet "演出" (2): refresh col F (and one col G) counts.
#  - Sheet "本地生活" (3): refresh col F count.
#  - Sheet "全部类型" (4): refresh the corresponding mirrored col F counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value  = 7661
$ws1.Range("F3").Value  = 7661
$ws1.Range("F5").Value  = 7842
$ws1.Range("F9").Value  = 6603
$ws1.Range("F10").Value = 3362
$ws1.Range("F12").Value = 3710
$ws1.Range("F13").Value = 42
$ws1.Range("F15").Value = 40
$ws1.Range("F20").Value = 26
$ws1.Range("F22").Value = 325
$ws1.Range("F28").Value = 1466
$ws1.Range("F30").Value = 53
$ws1.Range("F36").Value = 3643
$ws1.Range("F37").Value = 303
$ws1.Range("F39").Value = 42

# New row 46 — same look/feel (index cell bold+bordered) as the other rows.
$ws1.Range("A46").Value = 45
$ws1.Range("A45").Copy()
$ws1.Range("A46").PasteSpecial(-4122)   # xlPasteFormats - only copies A45's style onto A46

# Column B holds plain-text dates (not real Excel dates) in this sheet, so
# force text formatting before assignment to stop auto date-conversion, then
# drop back to the Normal style so no stray number-format sticks around.
$ws1.Range("B46").NumberFormat = "@"
$ws1.Range("B46").Value = "2024-10-03"
$ws1.Range("B46").Style = "Normal"

$ws1.Range("C46").Value = "北京·美漫超级英雄ONLY"
$ws1.Range("D46").Value = "小关路39号 北投购物公园"
$ws1.Range("E46").Value = "2024.10.03 09:30-10.04 17:00"
$ws1.Range("F46").Value = 0
$ws1.Range("G46").Value = "不可售"
$ws1.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=86413"
$ws1.Range("I46").Value = "//i0.hdslb.com/bfs/openplatform/202405/aPxrEklm1716799650037.png"

# ---------------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("G4").Value  = 380
$ws2.Range("F7").Value  = 40
$ws2.Range("F9").Value  = 102
$ws2.Range("F17").Value = 22

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F2").Value = 133

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (mirrors the three sheets above)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value  = 133
$ws4.Range("F5").Value  = 7661
$ws4.Range("F6").Value  = 7661
$ws4.Range("F7").Value  = 7842
$ws4.Range("F10").Value = 6603
$ws4.Range("F11").Value = 3362
$ws4.Range("F12").Value = 3710
$ws4.Range("F14").Value = 40
$ws4.Range("F18").Value = 40
$ws4.Range("F19").Value = 26
$ws4.Range("F22").Value = 325
$ws4.Range("F30").Value = 1466
$ws4.Range("F32").Value = 53
$ws4.Range("F38").Value = 3643
$ws4.Range("F39").Value = 303
$ws4.Range("F42").Value = 42
$ws4.Range("F45").Value = 22

Write-Output "edit applied"
